$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1945.9166
$ws.Range("I28").Value = 2087.7058
$ws.Range("J28").Value = 1601.5714
$ws.Range("K28").Value = 2087.7058
$ws.Range("L28").Value = 1601.5714
$ws.Range("M28").Value = -1602.7058
$ws.Range("N28").Value = -2571.5714
# Row 33
$ws.Range("H33").Value = 445.0909
$ws.Range("I33").Value = 459.4
$ws.Range("K33").Value = 459.4
$ws.Range("M33").Value = -230.4
# Row 88
$ws.Range("H88").Value = 2263.842
$ws.Range("I88").Value = 2089
$ws.Range("K88").Value = 2089
$ws.Range("M88").Value = -1683
# Row 91
$ws.Range("H91").Value = 2263.842
$ws.Range("I91").Value = 2089
$ws.Range("K91").Value = 2089
$ws.Range("M91").Value = -685
# Row 92
$ws.Range("H92").Value = 1008.75
$ws.Range("I92").Value = 1024
$ws.Range("K92").Value = 1024
$ws.Range("M92").Value = 224
# Row 107
$ws.Range("H107").Value = 1820.6666
$ws.Range("I107").Value = 1820.6666
$ws.Range("K107").Value = 1820.6666
$ws.Range("M107").Value = 99.33339999999998
# Row 127
$ws.Range("H127").Value = 1824.5
$ws.Range("I127").Value = 1656.8889
$ws.Range("K127").Value = 4970.6667
$ws.Range("M127").Value = -10.66669999999976
# Row 129
$ws.Range("H129").Value = 1716.762
$ws.Range("I129").Value = 760.8461
$ws.Range("K129").Value = 2282.5383
$ws.Range("M129").Value = 2717.4617
# Row 137
$ws.Range("H137").Value = 2775
$ws.Range("I137").Value = 2700
$ws.Range("K137").Value = 8100
$ws.Range("M137").Value = -5550

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 918.2
$ws.Range("I2").Value = 918.2
$ws.Range("K2").Value = 918.2
$ws.Range("M2").Value = -805.2
# Row 32
$ws.Range("H32").Value = 5450.724
$ws.Range("I32").Value = 5450.724
$ws.Range("K32").Value = 5450.724
$ws.Range("M32").Value = -5163.724
# Row 61
$ws.Range("H61").Value = 2195.5625
$ws.Range("I61").Value = 1594.5
$ws.Range("K61").Value = 1594.5
$ws.Range("M61").Value = -1382.5
# Row 97
$ws.Range("H97").Value = 1189.4546
$ws.Range("I97").Value = 803.7778
$ws.Range("K97").Value = 803.7778
$ws.Range("M97").Value = -307.7778
# Row 102
$ws.Range("H102").Value = 976.93335
$ws.Range("I102").Value = 976
$ws.Range("J102").Value = 990
$ws.Range("K102").Value = 976
$ws.Range("L102").Value = 990
$ws.Range("M102").Value = 646
$ws.Range("N102").Value = -4234
# Row 109
$ws.Range("H109").Value = 29999.5
$ws.Range("J109").Value = 29999.5
$ws.Range("L109").Value = 29999.5
$ws.Range("N109").Value = -32773.5
# Row 116
$ws.Range("H116").Value = 918.2
$ws.Range("I116").Value = 918.2
$ws.Range("K116").Value = 918.2
$ws.Range("M116").Value = 1375.8
# Row 132
$ws.Range("H132").Value = 1765
$ws.Range("I132").Value = 1735.1111
$ws.Range("J132").Value = 1899.5
$ws.Range("K132").Value = 5205.3333
$ws.Range("L132").Value = 5698.5
$ws.Range("M132").Value = -2675.3333
$ws.Range("N132").Value = -10758.5
# Row 136
$ws.Range("H136").Value = 2195.5625
$ws.Range("I136").Value = 1594.5
$ws.Range("K136").Value = 4783.5
$ws.Range("M136").Value = -2233.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 918.2
$ws.Range("I3").Value = 918.2
$ws.Range("K3").Value = 918.2
$ws.Range("M3").Value = -804.2
# Row 20
$ws.Range("H20").Value = 2648.3572
$ws.Range("I20").Value = 2570.9092
$ws.Range("K20").Value = 2570.9092
$ws.Range("M20").Value = -2323.9092
# Row 94
$ws.Range("H94").Value = 1213.8948
$ws.Range("I94").Value = 2021.5
$ws.Range("K94").Value = 2021.5
$ws.Range("M94").Value = -1570.5
# Row 107
$ws.Range("H107").Value = 1811.5834
$ws.Range("I107").Value = 1658.5555
$ws.Range("J107").Value = 2270.6667
$ws.Range("K107").Value = 1658.5555
$ws.Range("L107").Value = 2270.6667
$ws.Range("M107").Value = 261.4445000000001
$ws.Range("N107").Value = -6110.6667
# Row 134
$ws.Range("H134").Value = 1381.0714
$ws.Range("I134").Value = 1381.0714
$ws.Range("K134").Value = 4143.2142
$ws.Range("M134").Value = -1608.2142

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 205.5
$ws.Range("I2").Value = 205.5
$ws.Range("K2").Value = 205.5
$ws.Range("M2").Value = -92.5
# Row 31
$ws.Range("H31").Value = 2509
$ws.Range("I31").Value = 1499.6666
$ws.Range("J31").Value = 3518.3333
$ws.Range("K31").Value = 1499.6666
$ws.Range("L31").Value = 3518.3333
$ws.Range("M31").Value = -1204.6666
$ws.Range("N31").Value = -4108.3333
# Row 34
$ws.Range("H34").Value = 2509
$ws.Range("I34").Value = 1499.6666
$ws.Range("J34").Value = 3518.3333
$ws.Range("K34").Value = 1499.6666
$ws.Range("L34").Value = 3518.3333
$ws.Range("M34").Value = -1297.6666
$ws.Range("N34").Value = -3922.3333
# Row 44
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15884
# Row 58
$ws.Range("H58").Value = 5562.375
$ws.Range("I58").Value = 4699.8
$ws.Range("K58").Value = 4699.8
$ws.Range("M58").Value = -4496.8
# Row 105
$ws.Range("H105").Value = 2023.8
$ws.Range("I105").Value = 1954.75
$ws.Range("J105").Value = 2300
$ws.Range("K105").Value = 1954.75
$ws.Range("L105").Value = 2300
$ws.Range("M105").Value = -207.75
$ws.Range("N105").Value = -5794
# Row 132
$ws.Range("H132").Value = 2534
$ws.Range("I132").Value = 2570.818
$ws.Range("J132").Value = 2466.5
$ws.Range("K132").Value = 7712.454000000001
$ws.Range("L132").Value = 7399.5
$ws.Range("M132").Value = -5182.454000000001
$ws.Range("N132").Value = -12459.5
# Row 134
$ws.Range("H134").Value = 5642.8335
$ws.Range("I134").Value = 5642.8335
$ws.Range("K134").Value = 16928.5005
$ws.Range("M134").Value = -14393.5005
# Row 136
$ws.Range("H136").Value = 5562.375
$ws.Range("I136").Value = 4699.8
$ws.Range("K136").Value = 14099.4
$ws.Range("M136").Value = -11549.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 3561974.8
$ws.Range("I4").Value = 1818811.1
$ws.Range("J4").Value = 6002403.5
$ws.Range("K4").Value = 5456433.300000001
$ws.Range("L4").Value = 18007210.5
$ws.Range("M4").Value = -5456321.300000001
$ws.Range("N4").Value = -18007434.5
# Row 23
$ws.Range("H23").Value = 450.66666
$ws.Range("J23").Value = 479
$ws.Range("L23").Value = 1437
$ws.Range("N23").Value = -1907
# Row 32
$ws.Range("H32").Value = 5664.3335
$ws.Range("J32").Value = 7499
$ws.Range("L32").Value = 22497
$ws.Range("N32").Value = -23063
# Row 34
$ws.Range("H34").Value = 4209.55
$ws.Range("J34").Value = 4420.579
$ws.Range("L34").Value = 13261.737
$ws.Range("N34").Value = -13429.737
# Row 39
$ws.Range("H39").Value = 5980.85
$ws.Range("J39").Value = 6573.222
$ws.Range("L39").Value = 19719.666
$ws.Range("N39").Value = -20307.666
# Row 55
$ws.Range("H55").Value = 3670.077
$ws.Range("J55").Value = 4719.9
$ws.Range("L55").Value = 14159.7
$ws.Range("N55").Value = -14513.7
# Row 68
$ws.Range("H68").Value = 1224.25
$ws.Range("I68").Value = 900
$ws.Range("J68").Value = 1270.5714
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 3811.7142
$ws.Range("M68").Value = -1889
$ws.Range("N68").Value = -5433.7142
# Row 71
$ws.Range("H71").Value = 1224.25
$ws.Range("I71").Value = 900
$ws.Range("J71").Value = 1270.5714
$ws.Range("K71").Value = 8100
$ws.Range("L71").Value = 11435.1426
$ws.Range("M71").Value = -4044
$ws.Range("N71").Value = -19547.1426
# Row 122
$ws.Range("H122").Value = 1579.6
$ws.Range("J122").Value = 1524.5
$ws.Range("L122").Value = 13720.5
$ws.Range("N122").Value = -18620.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 1000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 8
$ws.Range("H8").Value = 1000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
# Row 132
$ws.Range("H132").Value = 1912.1111
$ws.Range("I132").Value = 1912.1111
$ws.Range("K132").Value = 5736.3333
$ws.Range("M132").Value = -3206.3333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 2659.2
$ws.Range("I122").Value = 2934.889
$ws.Range("K122").Value = 8804.667000000001
$ws.Range("M122").Value = -6354.667000000001
# Row 132
$ws.Range("H132").Value = 3379.05
$ws.Range("J132").Value = 2640
$ws.Range("L132").Value = 7920
$ws.Range("N132").Value = -12980

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2935.7144
$ws.Range("I81").Value = 3008.3333
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 6016.6666
$ws.Range("L81").Value = 5000
$ws.Range("M81").Value = -4955.6666
$ws.Range("N81").Value = -7122
# Row 84
$ws.Range("H84").Value = 2935.7144
$ws.Range("I84").Value = 3008.3333
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 30083.333
$ws.Range("L84").Value = 25000
$ws.Range("M84").Value = -24779.333
$ws.Range("N84").Value = -35608
# Row 104
$ws.Range("H104").Value = 14566.333
$ws.Range("J104").Value = 14566.333
$ws.Range("L104").Value = 14566.333
$ws.Range("N104").Value = -21554.333
# Row 126
$ws.Range("H126").Value = 1802.3572
$ws.Range("I126").Value = 2029.25
$ws.Range("K126").Value = 6087.75
$ws.Range("M126").Value = -3617.75
# Row 132
$ws.Range("H132").Value = 2918
$ws.Range("I132").Value = 2445.8667
$ws.Range("K132").Value = 7337.6001
$ws.Range("M132").Value = -4807.6001
